$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the taxon-record data held in row 22 and row 23 (columns A, B, D-I, Q, R).
# P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY are identical between the
# two rows already, so they are left untouched.

$cols = @("A","B","D","E","F","G","H","I","Q","R")

foreach ($col in $cols) {
    $addr22 = "$col" + "22"
    $addr23 = "$col" + "23"
    $val22 = $ws.Range($addr22).Value2
    $val23 = $ws.Range($addr23).Value2
    $ws.Range($addr22).Value = $val23
    $ws.Range($addr23).Value = $val22
}
